$wb = $excel.ActiveWorkbook

# --- Update the "Status" text from "Ready for handoff" to "In Translation" ---
# Sheet 1: Overview (columns E = zh-cn, F = de-de, rows 2-3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Sheet 2: zh-cn (column C = Status, rows 2-3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

# Sheet 3: de-de (column C = Status, rows 2-3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Narrow the Status columns to fit the shorter text ---
# Target stored column width is ~13.41 character-units; the closest width
# achievable through the ColumnWidth property (character units) is 12.5,
# which yields that stored width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
